$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reassert the header string (value text is unchanged, but the shared-string
# table was regenerated upstream with duplicate entries before this one).
$ws.Range("A1").Value = "HK_R_acc_SD"

# Updated, more reproducible match-percentage values for the data rows.
$values = @{
    2 = 100
    3 = 100
    4 = 100
    5 = 100
    6 = 100
    7 = 100
    8 = 99.8901098901099
    9 = 99.8901098901099
    10 = 99.8901098901099
    11 = 99.8901098901099
    12 = 99.8901098901099
    13 = 99.8901098901099
    14 = 99.8901098901099
    15 = 99.8901098901099
    16 = 99.8901098901099
    17 = 99.8901098901099
    18 = 99.8901098901099
    19 = 99.8901098901099
    20 = 100
    21 = 100
    22 = 100
    23 = 100
    24 = 100
    25 = 100
    26 = 99.8901098901099
    27 = 99.8901098901099
    28 = 99.8901098901099
    29 = 99.8901098901099
    30 = 99.8901098901099
    31 = 99.8901098901099
    32 = 100
    33 = 100
    34 = 100
    35 = 100
    36 = 100
    37 = 100
    38 = 100
    39 = 100
    40 = 100
    41 = 99.8901098901099
    42 = 99.8901098901099
    43 = 99.8901098901099
    44 = 99.8901098901099
    45 = 99.8901098901099
    46 = 99.8901098901099
    47 = 99.8901098901099
    48 = 99.8901098901099
    49 = 99.8901098901099
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
